$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (Q6)
$ws.Range("B7").Value = -0.03521679937216909
$ws.Range("C7").Value = 1.035216799372169
$ws.Range("D7").Value = 2.260394904707727
$ws.Range("E7").Value = 1.503460975452215
$ws.Range("F7").Value = 1.523224512798948
$ws.Range("G7").Value = 38

# Row 8 (Q7)
$ws.Range("B8").Value = -0.0694220599598384
$ws.Range("C8").Value = 0.9613139518517302
$ws.Range("D8").Value = 1.915446284370705
$ws.Range("E8").Value = 1.383996490013867
$ws.Range("F8").Value = 1.401320744993615
$ws.Range("G8").Value = 37

# Row 9 (Q8)
$ws.Range("B9").Value = 0.07999999999999999
$ws.Range("C9").Value = 1.11
$ws.Range("D9").Value = 2.875
$ws.Range("E9").Value = 1.695582495781317
$ws.Range("F9").Value = 1.737693571193845
$ws.Range("G9").Value = 20

# Row 10 (Q9)
$ws.Range("B10").Value = -0.353846153846154
$ws.Range("C10").Value = 1.353846153846154
$ws.Range("D10").Value = 3.603076923076923
$ws.Range("E10").Value = 1.898177263344212
$ws.Range("F10").Value = 1.94105443192413
$ws.Range("G10").Value = 13

# Row 11 (Q10 label, last row)
$ws.Range("B11").Value = -0.9199999999999999
$ws.Range("C11").Value = 1.16
$ws.Range("D11").Value = 2.124
$ws.Range("E11").Value = 1.457394936178934
$ws.Range("F11").Value = 1.26372465355393
$ws.Range("G11").Value = 5
